# Update the "Results" table on slide 14 to match the new experiment
# numbers: Noise Length's "Expert Results" value moves from 18 to 20, and
# the Threshold / Attack-Release rows no longer have values in the two
# "Expert Results" columns (replaced with a dash, matching the rest of
# the table's "no data" convention).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)
$shp = $s.Shapes.Item(3)
$tbl = $shp.Table

# Row 5 "Noise Length", column 3 (second "Expert Results" column): 18 -> 20
$tbl.Cell(5, 3).Shape.TextFrame.TextRange.Text = "20"

# Row 8 "Threshold", columns 2 & 3 ("Expert Results"): 0.03 -> -
$tbl.Cell(8, 2).Shape.TextFrame.TextRange.Text = "-"
$tbl.Cell(8, 3).Shape.TextFrame.TextRange.Text = "-"

# Row 9 "Attack/Release", columns 2 & 3 ("Expert Results"): 1.5 -> -
$tbl.Cell(9, 2).Shape.TextFrame.TextRange.Text = "-"
$tbl.Cell(9, 3).Shape.TextFrame.TextRange.Text = "-"
